# Apply crypto price/volume/name updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.798.63'
$ws.Range('E2').Value = '  +1.30%  '
$ws.Range('D3').Value = '2.524.70'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.87'
$ws.Range('E5').Value = '  +4.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.83'
$ws.Range('E6').Value = '  -1.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.579'
$ws.Range('E7').Value = '  -0.80%  '
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.529'
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.87'
$ws.Range('E10').Value = '  -1.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0811'
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.56'
$ws.Range('E12').Value = '  +1.66%  '
$ws.Range('E13').Value = '  -1.97%  '
$ws.Range('D14').Value = '2.910.21'
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.540.31'
$ws.Range('E15').Value = '  +0.66%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.22'
$ws.Range('E16').Value = '  +1.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.847'
$ws.Range('E17').Value = '  -0.83%  '
$ws.Range('D18').Value = '42.852.32'
$ws.Range('E18').Value = '  +1.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.96'
$ws.Range('E19').Value = '  +1.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.68'
$ws.Range('E20').Value = '  +4.05%  '
$ws.Range('D21').Value = '0.0₃0963'
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.79'
$ws.Range('E22').Value = '  -1.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '250.83'
$ws.Range('E23').Value = '  +0.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.94'
$ws.Range('E24').Value = '  +2.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.02'
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.73'
$ws.Range('E26').Value = '  -0.80%  '
$ws.Range('E27').Value = '  -0.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.43'
$ws.Range('E28').Value = '  +4.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '40.34'
$ws.Range('E29').Value = '  +6.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.25'
$ws.Range('E30').Value = '  +0.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.00'
$ws.Range('E31').Value = '  +1.62%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '156.15'
$ws.Range('E32').Value = '  +0.93%  '
$ws.Range('E33').Value = '  +3.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.99'
$ws.Range('E34').Value = '  +2.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.29'
$ws.Range('E35').Value = '  -0.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0788'
$ws.Range('E36').Value = '  +0.63%  '
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('E38').Value = '  -2.51%  '
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.60'
$ws.Range('E40').Value = '  -1.93%  '
$ws.Range('E41').Value = '  +14.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0304'
$ws.Range('E42').Value = '  +1.96%  '
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('E44').Value = '  -1.41%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.29'
$ws.Range('E45').Value = '  -2.08%  '
$ws.Range('D46').Value = '2.019.58'
$ws.Range('E46').Value = '  -0.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.68'
$ws.Range('E47').Value = '  +1.54%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.76'
$ws.Range('E48').Value = '  -1.91%  '
$ws.Range('D49').Value = '2.765.29'
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.64'
$ws.Range('E50').Value = '  +2.98%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '102.59'
$ws.Range('E51').Value = '  +1.40%  '
